# The underlying data table (rows 2-285, columns A:R) records weekly price
# observations for "Cebollín" at Femacal de La Calera. This commit adds a new
# week of observations: two new rows are inserted right after the existing
# row for 2021-07-12 (old row 197), pushing that row and everything below it
# down by two, and the two freshly inserted rows are populated with the new
# week's "Primera" / "Segunda" quality records (2021-09-27).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the existing row 197 - this shifts the
# old rows 197..285 down to 199..287 (and therefore also pushes the old
# 283/284/285 tail down to 285/286/287), matching the new A1:R287 extent.
$ws.Rows.Item(197).Resize(2).Insert()

# Populate the first newly-inserted row (197) - "Primera" quality for the
# new 2021-09-27 observation.
$ws.Cells.Item(197, 1).Value = 3
$ws.Cells.Item(197, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(197, 3).Value = "Coquimbo"
$ws.Cells.Item(197, 4).Value = 44466
$ws.Cells.Item(197, 5).Value = 5
$ws.Cells.Item(197, 6).Value = 100112037
$ws.Cells.Item(197, 7).Value = "Cebollín"
$ws.Cells.Item(197, 8).Value = "Sin especificar"
$ws.Cells.Item(197, 9).Value = "Primera"
$ws.Cells.Item(197, 10).Value = 270
$ws.Cells.Item(197, 11).Value = 3300
$ws.Cells.Item(197, 12).Value = 3500
$ws.Cells.Item(197, 13).Value = 3381
$ws.Cells.Item(197, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(197, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(197, 16).Value = 94
$ws.Cells.Item(197, 17).Value = 36
$ws.Cells.Item(197, 18).Value = "Hortaliza"

# Populate the second newly-inserted row (198) - "Segunda" quality for the
# same new 2021-09-27 observation.
$ws.Cells.Item(198, 1).Value = 3
$ws.Cells.Item(198, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(198, 3).Value = "Coquimbo"
$ws.Cells.Item(198, 4).Value = 44466
$ws.Cells.Item(198, 5).Value = 5
$ws.Cells.Item(198, 6).Value = 100112037
$ws.Cells.Item(198, 7).Value = "Cebollín"
$ws.Cells.Item(198, 8).Value = "Sin especificar"
$ws.Cells.Item(198, 9).Value = "Segunda"
$ws.Cells.Item(198, 10).Value = 120
$ws.Cells.Item(198, 11).Value = 2500
$ws.Cells.Item(198, 12).Value = 2500
$ws.Cells.Item(198, 13).Value = 2500
$ws.Cells.Item(198, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(198, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(198, 16).Value = 69
$ws.Cells.Item(198, 17).Value = 36
$ws.Cells.Item(198, 18).Value = "Hortaliza"
